$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "39÷8="
$t.Cell(1, 2).Range.Text = "53÷9="
$t.Cell(1, 3).Range.Text = "63÷8="
$t.Cell(1, 4).Range.Text = "44÷8="
$t.Cell(1, 5).Range.Text = "51÷4="
$t.Cell(5, 1).Range.Text = "58÷8="
$t.Cell(5, 2).Range.Text = "94÷7="
$t.Cell(5, 3).Range.Text = "40÷6="
$t.Cell(5, 4).Range.Text = "96÷4="
$t.Cell(5, 5).Range.Text = "23÷9="
$t.Cell(9, 1).Range.Text = "98÷7="
$t.Cell(9, 2).Range.Text = "80÷5="
$t.Cell(9, 3).Range.Text = "56÷8="
$t.Cell(9, 4).Range.Text = "15÷2="
$t.Cell(9, 5).Range.Text = "10÷6="
$t.Cell(13, 1).Range.Text = "57÷9="
$t.Cell(13, 2).Range.Text = "21÷5="
$t.Cell(13, 3).Range.Text = "60÷3="
$t.Cell(13, 4).Range.Text = "66÷2="
$t.Cell(13, 5).Range.Text = "83÷9="
$t.Cell(17, 1).Range.Text = "56÷9="
$t.Cell(17, 2).Range.Text = "14÷6="
$t.Cell(17, 3).Range.Text = "98÷2="
$t.Cell(17, 4).Range.Text = "33÷8="
$t.Cell(17, 5).Range.Text = "96÷8="
